# project-general.xlsx: insert a blank spacer row directly above each of the
# two "Total" rows (the "with project" total and the "without project" total).
# Excel's EntireRow Insert pushes everything below down by one row, copies the
# formatting from the adjoining row, and auto-updates formulas / data
# validation ranges that reference the shifted cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blank row above the first "Total" row (old row 14 -> new row 15).
$ws.Rows(14).Insert()

# Blank row above the second "Total" row. After the first insert shifted
# everything down by one, that row now lives at row 30.
$ws.Rows(30).Insert()

# Leave the sheet's selection on the newly inserted second blank row, as it
# was left positioned in the saved workbook.
$ws.Range("A30:XFD30").Select()
